# ComputerPartsData.xlsx edit
# - Strip the leading "$" currency sign from every Price value in column C
#   (and correct a handful of stale prices that had drifted from the live site).
# - Fix WebFunctions-style append bug: rows 47, 54, 56, 63, 64 were built from a
#   "search_list" whose matched entries were not removed after use, so later scrapes
#   appended into the wrong / already-consumed row and left every column after Brand
#   filled in as "DNE". Row 65 had in turn absorbed what should have been row 64's
#   real data. We reproduce that corrected (still bug-affected per upstream data) state
#   here: the leftover price leaks into column A and the remaining columns become "DNE".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Strip "$" from Price column (column C) ---
Set-TextValue $ws.Range("C2") "159"
Set-TextValue $ws.Range("C3") "159"
Set-TextValue $ws.Range("C4") "589"
Set-TextValue $ws.Range("C5") "199"
Set-TextValue $ws.Range("C6") "317"
Set-TextValue $ws.Range("C7") "249"
Set-TextValue $ws.Range("C8") "177"
Set-TextValue $ws.Range("C9") "358"
Set-TextValue $ws.Range("C10") "279"
Set-TextValue $ws.Range("C11") "137"
Set-TextValue $ws.Range("C12") "391"
Set-TextValue $ws.Range("C13") "288"
Set-TextValue $ws.Range("C14") "399"
Set-TextValue $ws.Range("C15") "148"
Set-TextValue $ws.Range("C16") "124"
Set-TextValue $ws.Range("C17") "161"
Set-TextValue $ws.Range("C18") "109"
Set-TextValue $ws.Range("C19") "499"
Set-TextValue $ws.Range("C20") "193"
Set-TextValue $ws.Range("C21") "249"
Set-TextValue $ws.Range("C22") "203"
Set-TextValue $ws.Range("C23") "448"
Set-TextValue $ws.Range("C24") "349"
Set-TextValue $ws.Range("C25") "154"
Set-TextValue $ws.Range("C26") "328"
Set-TextValue $ws.Range("C27") "112"
Set-TextValue $ws.Range("C28") "129"
Set-TextValue $ws.Range("C29") "312"
Set-TextValue $ws.Range("C30") "734"
Set-TextValue $ws.Range("C31") "224"
Set-TextValue $ws.Range("C32") "298"
Set-TextValue $ws.Range("C33") "379"
Set-TextValue $ws.Range("C34") "159"
Set-TextValue $ws.Range("C35") "389"
Set-TextValue $ws.Range("C36") "377"
Set-TextValue $ws.Range("C37") "259"
Set-TextValue $ws.Range("C38") "149"
Set-TextValue $ws.Range("C39") "279"
Set-TextValue $ws.Range("C40") "289"
Set-TextValue $ws.Range("C41") "312"
Set-TextValue $ws.Range("C42") "342"
Set-TextValue $ws.Range("C43") "313"
Set-TextValue $ws.Range("C44") "146"
Set-TextValue $ws.Range("C45") "229"
Set-TextValue $ws.Range("C46") "379"
Set-TextValue $ws.Range("C48") "249"
Set-TextValue $ws.Range("C49") "209"
Set-TextValue $ws.Range("C50") "179"
Set-TextValue $ws.Range("C51") "224"
Set-TextValue $ws.Range("C52") "110"
Set-TextValue $ws.Range("C53") "549"
Set-TextValue $ws.Range("C55") "202"
Set-TextValue $ws.Range("C57") "359"
Set-TextValue $ws.Range("C58") "288"
Set-TextValue $ws.Range("C59") "163"
Set-TextValue $ws.Range("C60") "76"
Set-TextValue $ws.Range("C61") "246"
Set-TextValue $ws.Range("C62") "139"
Set-TextValue $ws.Range("C66") "77"
Set-TextValue $ws.Range("C68") "149"
Set-TextValue $ws.Range("C69") "142"
Set-TextValue $ws.Range("C70") "459"

# --- Rows whose data shifted due to the WebFunctions append bug ---
# Row 47
Set-TextValue $ws.Range("A47") "428"
Set-TextValue $ws.Range("B47") "DNE"
Set-TextValue $ws.Range("C47") "DNE"
Set-TextValue $ws.Range("D47") "DNE"
Set-TextValue $ws.Range("E47") "DNE"
Set-TextValue $ws.Range("F47") "DNE"
# Row 54
Set-TextValue $ws.Range("A54") "219"
Set-TextValue $ws.Range("B54") "DNE"
Set-TextValue $ws.Range("C54") "DNE"
Set-TextValue $ws.Range("D54") "DNE"
Set-TextValue $ws.Range("E54") "DNE"
Set-TextValue $ws.Range("F54") "DNE"
# Row 56
Set-TextValue $ws.Range("A56") "199"
Set-TextValue $ws.Range("B56") "DNE"
Set-TextValue $ws.Range("C56") "DNE"
Set-TextValue $ws.Range("D56") "DNE"
Set-TextValue $ws.Range("E56") "DNE"
Set-TextValue $ws.Range("F56") "DNE"
# Row 63
Set-TextValue $ws.Range("A63") "3,374"
Set-TextValue $ws.Range("B63") "DNE"
Set-TextValue $ws.Range("C63") "DNE"
Set-TextValue $ws.Range("D63") "DNE"
Set-TextValue $ws.Range("E63") "DNE"
Set-TextValue $ws.Range("F63") "DNE"
# Row 64
Set-TextValue $ws.Range("A64") "7,490"
Set-TextValue $ws.Range("B64") "DNE"
Set-TextValue $ws.Range("C64") "DNE"
Set-TextValue $ws.Range("D64") "DNE"
Set-TextValue $ws.Range("E64") "DNE"
Set-TextValue $ws.Range("F64") "DNE"
Set-TextValue $ws.Range("G64") "DNE"
Set-TextValue $ws.Range("H64") "DNE"
# Row 65
Set-TextValue $ws.Range("A65") "AMD"
Set-TextValue $ws.Range("B65") "Ryzen Threadripper 3960X"
Set-TextValue $ws.Range("C65") "1,532"
Set-TextValue $ws.Range("D65") "Socket sTRX4"
Set-TextValue $ws.Range("E65") "24-Core"
Set-TextValue $ws.Range("F65") "48"
Set-TextValue $ws.Range("G65") "3.8 GHz"
Set-TextValue $ws.Range("H65") "Up to 4.5 GHz"
# Row 67
Set-TextValue $ws.Range("A67") "99"
